# TC05_Search_product_in_Catalog.xlsx
#
# "Reexecution of failed testcases logic implementation"
#
# A new test step is inserted as the new row 3 on the main sheet:
#   Keyword=CLICK_PRE_ENTERTEXT, Object=SearchBoxHomePage, ObjectType=CSS
# All the former rows 3-7 shift down to rows 4-8 (their content is
# untouched). The sheet's used range grows from A1:E7 to A1:E8 and the
# active selection moves to the new row's B:D cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the existing data rows (old rows 3-7) down by one to make room for
# the new step - mirrors Excel's "Insert Sheet Rows" at row 3.
$ws.Rows.Item(3).Insert()

# Fill in the new step's data (A3 / E3 stay blank, matching the other rows
# whose Data_descriptor column is empty).
$ws.Range("B3").Value2 = "CLICK_PRE_ENTERTEXT"
$ws.Range("C3").Value2 = "SearchBoxHomePage"
$ws.Range("D3").Value2 = "CSS"

# Give the new row the same bordered look as the surrounding data rows.
$ws.Range("A3:E3").Borders.LineStyle = 1

# Match the workbook's saved selection (B3:D3, active cell B3).
$ws.Range("B3:D3").Select()
